$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column D so that the old "Terms Typically
# Offered" column (D) shifts to G, making room for the new
# Corequisites / Concurrent / Recommended columns.
$ws.Range("D1:F1").EntireColumn.Insert()

# New header row
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Default all new cells (D2:F20) to "NA"
$ws.Range("D2:F20").Value = "NA"

# A handful of prerequisite cells used a non-breaking space between "WLC"
# and the course number; normalize these to a regular space to match the
# updated data.
$ws.Range("C3").Value = "WLC 101 or consent of instructor."
$ws.Range("C4").Value = "WLC 102 or consent of instructor."
$ws.Range("C6").Value = "WLC 103."
$ws.Range("C7").Value = "WLC 201."
$ws.Range("C19").Value = "WLC 360, advanced composition in primary and/or secondary language, senior standing and consent of instructor."

# Row 12 - WLC 318: move "Corequisite: WLC 310." out of prerequisites
# column into the new Corequisites column.
$ws.Range("C12").Value = "Limited to Valladolid, Spain Fall program."
$ws.Range("D12").Value = "WLC 310."
$ws.Range("G12").Value = "F, SU "

# Row 13 - WLC 360: move "Recommended: SPAN 233 and SPAN 301." out of
# prerequisites column into the new Recommended column.
$ws.Range("C13").Value = "Junior standing; Modern Languages and Literatures major or Spanish major."
$ws.Range("F13").Value = "SPAN 233 and SPAN 301."
$ws.Range("G13").Value = "W "

# Row 17 - WLC 425: move "Corequisite: Concurrent enrollment in EDUC 469
# or EDUC 479." out of prerequisites column into the new Corequisites
# column.
$ws.Range("C17").Value = "Admission to the Single Subject Credential Program in World Languages."
$ws.Range("D17").Value = "Concurrent enrollment in EDUC 469 or EDUC 479."
$ws.Range("G17").Value = "W, SP "
